$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.965.97"
$ws.Range("E2").Value = "  -0.27%  "
Set-TextValue $ws.Range("D3") "1.636.06"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "214.63"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.59%  "
Set-TextValue $ws.Range("D9") "0.0624"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("E10").Value = "  -5.75%  "
$ws.Range("E11").Value = "  -0.58%  "
Set-TextValue $ws.Range("D12") "1.864.30"
$ws.Range("E12").Value = "  -0.54%  "
Set-TextValue $ws.Range("D13") "1.638.63"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("E15").Value = "  -2.27%  "
Set-TextValue $ws.Range("D16") "25.986.54"
$ws.Range("E17").Value = "  -2.56%  "
Set-TextValue $ws.Range("D18") "61.77"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").Value = "  -0.11%  "
Set-TextValue $ws.Range("D20") "191.11"
$ws.Range("E20").Value = "  -1.08%  "
Set-TextValue $ws.Range("D22") "9.65"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  +1.23%  "
Set-TextValue $ws.Range("D25") "143.53"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D26") "1.78"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D27") "1.01"
$ws.Range("E27").Value = "  -0.10%  "
Set-TextValue $ws.Range("D28") "6.82"
$ws.Range("E28").Value = "  -1.42%  "
Set-TextValue $ws.Range("D29") "15.29"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("E30").Value = "  -1.40%  "
Set-TextValue $ws.Range("D31") "0.0484"
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -2.57%  "
Set-TextValue $ws.Range("D36") "1.138.69"
$ws.Range("E36").Value = "  +0.81%  "
Set-TextValue $ws.Range("D37") "0.867"
$ws.Range("E37").Value = "  -4.32%  "
Set-TextValue $ws.Range("D38") "2.44"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("E40").Value = "  -1.15%  "
Set-TextValue $ws.Range("D41") "98.57"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("E42").Value = "  -2.43%  "
Set-TextValue $ws.Range("D43") "5.25"
$ws.Range("E43").Value = "  -4.90%  "
Set-TextValue $ws.Range("D44") "1.773.63"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -1.54%  "
Set-TextValue $ws.Range("D46") "55.33"
$ws.Range("E46").Value = "  -2.33%  "
Set-TextValue $ws.Range("D47") "0.0528"
$ws.Range("E47").Value = "  -0.28%  "
Set-TextValue $ws.Range("D48") "1.50"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("E51").Value = "  -0.04%  "
